$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(34, 'Cortider', 1, 'Cortider 10gm Cream', '10 gm', 1161, 30464),
    @(53, 'Esoral', 2, 'Esoral 20mg Capsule 50''s', '50''s', 0, 0),
    @(53, 'Esoral', 3, 'Esoral 20mg Tablet  80''s', '80''s', 0, 0),
    @(53, 'Esoral', 4, 'Esoral Injection & Capsule 20', 'Bundle', 0, 0),
    @(53, 'Esoral', 5, 'Esoral 20mg Capsule Container 30''s', '30''s', 0, 0),
    @(53, 'Esoral', 6, 'Esoral 20mg Tablet - 50''s', '50''s', 0, 0),
    @(53, 'Esoral', 7, 'Esoral 40mg Tablet', '20''s', 0, 0),
    @(53, 'Esoral', 8, 'Esoral 40mg EC Tablet - 42''s', '42''s', 0, 0),
    @(53, 'Esoral', 9, 'Esoral Mups 40 Tablet', '30''s', 5162, 1625204),
    @(53, 'Esoral', 10, 'Esoral 40mg Tablet - 30''s', '30''s', 4221, 759400),
    @(53, 'Esoral', 11, 'Esoral 40mg IV Injection', '1''s', 51100, 1955086),
    @(53, 'Esoral', 12, 'Esoral 20mg Capsule 100''s', '100''s', 21275, 11163843),
    @(53, 'Esoral', 13, 'Esoral 40mg Capsule - 30''s', '30''s', 13082, 2942010),
    @(53, 'Esoral', 14, 'Esoral 20mg Tablet 112''s', '112''s', 22459, 9428063),
    @(53, 'Esoral', 15, 'Esoral MUPS 20mg Tab 50s', '50''s', 40322, 15113492),
    @(53, 'Esoral', 16, 'Esoral 20mg Tablet', '20''s', 0, 0),
    @(53, 'Esoral', 17, 'Esoral 20mg  Tablet 100''s', '100''s', 0, 0),
    @(53, 'Esoral', 18, 'Esoral 20mg Capsule 60''s', '60''s', 0, 0),
    @(53, 'Esoral', 19, 'Esoral Injection & MUPS 20', 'Bundle', 0, 0),
    @(56, 'ETRIAM', 20, 'Etriam 10gm Cream', '1''s', 36775, 1102882),
    @(57, 'Facid', 21, 'Facid 10gm BT Cream', '10 gm', 10467, 1333914),
    @(57, 'Facid', 22, 'Facid 15gm Cream', '15gm', 8243, 556155),
    @(57, 'Facid', 23, 'Facid 15gm Ointment', '15gm', 1786, 120501),
    @(57, 'Facid', 24, 'Facid 250mg Tablet 10''s', '10''s', 4232, 2062084),
    @(57, 'Facid', 25, 'Facid HC Cream', '10 gm', 10164, 952366),
    @(64, 'Flucoder', 26, 'Flucoder 150mg Capsule', '10 ''s', 693, 114289),
    @(64, 'Flucoder', 27, 'Flucoder 200mg Capsule - 10''s', '10''s', 260, 48726),
    @(64, 'Flucoder', 28, 'Flucoder 35ml Dry Suspension', '35ml', 1301, 76069),
    @(64, 'Flucoder', 29, 'Flucoder 50mg Capsule', '30 ''s', 2018, 367598),
    @(75, 'Hairgrow', 30, 'Hairgrow 2% Topical Solution 60ml', '1''s', 0, 0),
    @(75, 'Hairgrow', 31, 'Hairgrow 5% Topical Solution 60ml', '1''s', 80, 35982),
    @(85, 'Kezona', 32, 'Kezona 200mg Tablet 30''s', '30''s', 1533, 310279),
    @(92, 'Licnil', 33, 'Licnil 117gm Lotion', '1''s', 297, 46753),
    @(92, 'Licnil', 34, 'Licnil 60gm Lotion', '1''s', 1184, 115380),
    @(96, 'Losectil', 35, 'Losectil 40mg Capsule (24''s)', '24 ''s', 0, 0),
    @(96, 'Losectil', 36, 'Losectil 20mg Powder for Oral Suspension - 30''s', '30''s', 0, 0),
    @(96, 'Losectil', 37, 'Losectil 10mg Capsule', '48 ''s', 0, 0),
    @(96, 'Losectil', 38, 'Losectil DR Tablet ', '60 ''s', 0, 0),
    @(96, 'Losectil', 39, 'Losectil 40mg Powder for Oral Suspension', '20''s', 0, 0),
    @(96, 'Losectil', 40, 'Losectil 20mg Powder for Oral Suspension', '20''s', 0, 0),
    @(96, 'Losectil', 41, 'Losectil 40mg Powder for Oral Suspension - 30''s', '30''s', 3661, 823322),
    @(96, 'Losectil', 42, 'Losectil Mups 40mg Tablet', '24''s', 776, 139610),
    @(96, 'Losectil', 43, 'Losectil 20mg Capsule 120''s', '120''s', 166775, 75010391),
    @(96, 'Losectil', 44, 'Losectil MUPS 20mg Tablet', '50''s', 4330, 1136062),
    @(96, 'Losectil', 45, 'Losectil V 40mg Capsule', '60''s', 808, 290734),
    @(96, 'Losectil', 46, 'Losectil 40mg IV Injection', '1''s', 18807, 719555),
    @(96, 'Losectil', 47, 'Losectil 40mg Capsule - 48''s', '48''s', 0, 0),
    @(96, 'Losectil', 48, 'Losectil 20mg PFS 50''s', '50''s', 20492, 4608445),
    @(96, 'Losectil', 49, 'Losectil V 20mg Capsule', '100''s', 7170, 2687387),
    @(96, 'Losectil', 50, 'Losectil 20mg Capsule (100''s)', '100 ''s', 0, 0),
    @(96, 'Losectil', 51, 'Losectil 20mg Capsule 500s', '500''s', 0, 0),
    @(96, 'Losectil', 52, 'Losectil 10mg Capsule 60''s', '60''s', 13138, 1181894),
    @(96, 'Losectil', 53, 'Losectil 40mg Capsule 60''s', '60''s', 7232, 2927513),
    @(99, 'Lulizol', 54, 'Lulizol 20gm Cream', '1''s', 14624, 1973216),
    @(99, 'Lulizol', 55, 'Lulizol 10gm Cream', '1''s', 47610, 3568845),
    @(112, 'Mupiron', 56, 'Mupiron 20gm Ointment', '1''s', 225, 37107),
    @(112, 'Mupiron', 57, 'Mupiron 10gm Ointment', '10gm', 3078, 323036),
    @(113, 'Mycofin', 58, 'Mycofin 250mg Tablet', '10 ''s', 8371, 2510044),
    @(113, 'Mycofin', 59, 'Mycofin 5gm Cream', '5gm', 7345, 275290),
    @(121, 'Nospot', 60, 'Nospot 30gm Cream', '1''s', 2852, 427600),
    @(127, 'Orogurd', 61, 'Orogurd 15gm Oral Gel', '15gm', 6269, 281979),
    @(133, 'Panoral', 62, 'Panoral 20mg Tablet', '40 ''s', 824, 123534),
    @(133, 'Panoral', 63, 'Panoral 40mg Tablet', '40''s', 152, 31904),
    @(134, 'Perosa', 64, 'Perosa 30gm Cream', '30 gm', 4425, 165849),
    @(143, 'Rabifast', 65, 'Rabifast 20mg Tablet - 50''s', '50''s', 0, 0),
    @(143, 'Rabifast', 66, 'Rabifast 20 mg Tablet', '70''s', 7445, 1953344),
    @(165, 'Softi', 67, 'Softi Ointment 15gm', '15gm', 0, 0),
    @(165, 'Softi', 68, 'Softi Ointment 50gm', '50gm', 1948, 131431),
    @(170, 'Sorex', 69, 'Sorex 5gm Oral Paste', '5gm', 581, 32663),
    @(186, 'Topibet', 70, 'Topibet CL 10gm Ointment', '10gm', 9257, 208189),
    @(187, 'Topiclo', 71, 'Topiclo S 6% 10gm Ointment', '1''s', 0, 0),
    @(187, 'Topiclo', 72, 'Topiclo 10gm Ointment', '10gm', 996, 37330),
    @(187, 'Topiclo', 73, 'Topiclo 10gm Cream', '10gm', 728, 24555),
    @(187, 'Topiclo', 74, 'Topiclo S Ointment 1''s 20g', '1''s', 2031, 182708),
    @(187, 'Topiclo', 75, 'Topiclo S 10gm Ointment', '10gm', 21014, 1102604),
    @(187, 'Topiclo', 76, 'Topiclo NN 15gm Ointment', '1''s', 4146, 202034)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}